$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.585.38'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.664.68'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.10'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.89%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.13%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.22%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.89'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.79%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.57%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0879'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.902.28'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.692.74'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.13'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.559'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.44'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.583.22'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '240.98'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0${sub3}0730"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.65'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.67%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.49'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -3.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.41'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.21'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.41'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.96%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.28%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.59%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.462.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.12'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.58'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.01%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.925'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.08%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0173'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.575'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.77%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.03'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.51%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.61'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.22'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.07%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.808.60'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.788'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.74'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.08'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.04%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.52%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.88'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.91%  '
